$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed date) column C for rows 2-6
# from serial date 45175 (2023-09-06) to 45183 (2023-09-14)
$ws.Range("C2:C6").Value = 45183
